$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-07 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-08 Sunday", 2)

$d.Content.Find.Execute("241×4=", $true, $false, $false, $false, $false, $true, 1, $false, "964×6=", 2)
$d.Content.Find.Execute("318×2=", $true, $false, $false, $false, $false, $true, 1, $false, "872×6=", 2)
$d.Content.Find.Execute("225×5=", $true, $false, $false, $false, $false, $true, 1, $false, "880×9=", 2)
$d.Content.Find.Execute("451×2=", $true, $false, $false, $false, $false, $true, 1, $false, "143×7=", 2)
$d.Content.Find.Execute("822×6=", $true, $false, $false, $false, $false, $true, 1, $false, "639×2=", 2)
$d.Content.Find.Execute("113×2=", $true, $false, $false, $false, $false, $true, 1, $false, "318×8=", 2)
$d.Content.Find.Execute("712×2=", $true, $false, $false, $false, $false, $true, 1, $false, "792×7=", 2)
$d.Content.Find.Execute("255×5=", $true, $false, $false, $false, $false, $true, 1, $false, "969×9=", 2)
$d.Content.Find.Execute("470×6=", $true, $false, $false, $false, $false, $true, 1, $false, "633×9=", 2)
$d.Content.Find.Execute("571×7=", $true, $false, $false, $false, $false, $true, 1, $false, "580×9=", 2)
$d.Content.Find.Execute("517×6=", $true, $false, $false, $false, $false, $true, 1, $false, "683×7=", 2)
$d.Content.Find.Execute("773×5=", $true, $false, $false, $false, $false, $true, 1, $false, "521×2=", 2)
$d.Content.Find.Execute("141×5=", $true, $false, $false, $false, $false, $true, 1, $false, "869×6=", 2)
$d.Content.Find.Execute("745×4=", $true, $false, $false, $false, $false, $true, 1, $false, "260×4=", 2)
$d.Content.Find.Execute("390×5=", $true, $false, $false, $false, $false, $true, 1, $false, "616×2=", 2)
$d.Content.Find.Execute("307×3=", $true, $false, $false, $false, $false, $true, 1, $false, "793×9=", 2)
$d.Content.Find.Execute("343×2=", $true, $false, $false, $false, $false, $true, 1, $false, "950×7=", 2)
$d.Content.Find.Execute("954×2=", $true, $false, $false, $false, $false, $true, 1, $false, "629×6=", 2)
$d.Content.Find.Execute("465×5=", $true, $false, $false, $false, $false, $true, 1, $false, "501×5=", 2)
$d.Content.Find.Execute("883×7=", $true, $false, $false, $false, $false, $true, 1, $false, "986×7=", 2)
$d.Content.Find.Execute("643×4=", $true, $false, $false, $false, $false, $true, 1, $false, "595×5=", 2)
$d.Content.Find.Execute("991×6=", $true, $false, $false, $false, $false, $true, 1, $false, "734×4=", 2)
$d.Content.Find.Execute("829×7=", $true, $false, $false, $false, $false, $true, 1, $false, "797×6=", 2)
$d.Content.Find.Execute("475×8=", $true, $false, $false, $false, $false, $true, 1, $false, "791×6=", 2)
$d.Content.Find.Execute("828×8=", $true, $false, $false, $false, $false, $true, 1, $false, "277×8=", 2)
